$wb = $excel.ActiveWorkbook

# Rename the sheet from "SCP 1.0.5" to "SCP 1.1.0"
$ws = $wb.ActiveSheet
$ws.Name = "SCP 1.1.0"

# Update the selected/active cell on the sheet from B42 to C13
$ws.Range("C13").Select()
